# Applies the "automatic update of files" edit:
#  - Column C ("Förändrad") for data rows 2..32 is bumped from serial 45610 to 45611
#    (i.e. incremented by one day)
#  - Rows 30 and 31 have their "Beteckning" (A) and "Area (ha)" (G) values swapped:
#      A30: A 46082-2024 -> A 46085-2024 , G30: 1   -> 1.9
#      A31: A 46085-2024 -> A 46082-2024 , G31: 1.9 -> 1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment the "Förändrad" date (column C) for rows 2 through 32 by one day.
# Use Value2 (returns a plain numeric Double for date-formatted cells) so the
# arithmetic is numeric and the cell stays stored as a number, not a string.
for ($row = 2; $row -le 32; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = $current + 1
    }
}

# Swap the Beteckning / Area values between row 30 and row 31.
$ws.Range("A30").Value = "A 46085-2024"
$ws.Range("G30").Value2 = 1.9

$ws.Range("A31").Value = "A 46082-2024"
$ws.Range("G31").Value2 = 1
